$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 3.8
$ws.Range("J2").Value = 2.4
$ws.Range("K2").Value = 2.3
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 4.33
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.15
$ws.Range("S2").Value = 1.33
$ws.Range("T2").Value = 3.25
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 15
$ws.Range("AH2").Value = 13
$ws.Range("AI2").Value = 21
$ws.Range("AO2").Value = 9.5
$ws.Range("AT2").Value = 3.25
$ws.Range("G9").Value = 2.05
$ws.Range("H9").Value = 3.25
$ws.Range("J9").Value = 2.75
$ws.Range("K9").Value = 2
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 7.5
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 2.75
$ws.Range("Q9").Value = 2.3
$ws.Range("R9").Value = 1.6
$ws.Range("S9").Value = 1.5
$ws.Range("T9").Value = 2.5
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.75
$ws.Range("W9").Value = 6.5
$ws.Range("X9").Value = 9
$ws.Range("Y9").Value = 9
$ws.Range("AA9").Value = 19
$ws.Range("AB9").Value = 34
$ws.Range("AC9").Value = 7.5
$ws.Range("AE9").Value = 17
$ws.Range("AG9").Value = 401
$ws.Range("AH9").Value = 9
$ws.Range("AR9").Value = 67
$ws.Range("AS9").Value = 201
$ws.Range("AT9").Value = 2.5
$ws.Range("AY9").Value = 21
$ws.Range("BB9").Value = 101
$ws.Range("BC9").Value = 301
